$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.41"
$ws.Range("E2").Value = "'-0.85%"
$ws.Range("D3").Value = "'37.66"
$ws.Range("E3").Value = "'7.49%"
$ws.Range("D4").Value = "'4.969"
$ws.Range("E4").Value = "'-3.66%"
$ws.Range("D5").Value = "'0.07737"
$ws.Range("E5").Value = "'-0.48%"
$ws.Range("D6").Value = "'2.191"
$ws.Range("E6").Value = "'-7.90%"
$ws.Range("D7").Value = "'7.998"
$ws.Range("E7").Value = "'-0.44%"
$ws.Range("D8").Value = "'3.991"
$ws.Range("E8").Value = "'1.37%"
$ws.Range("D9").Value = "'0.9124"
$ws.Range("E9").Value = "'-2.16%"
$ws.Range("D10").Value = "'0.09375"
$ws.Range("E10").Value = "'-6.70%"
$ws.Range("D11").Value = "'0.1796"
$ws.Range("E11").Value = "'-0.03%"
$ws.Range("D12").Value = "'0.08443"
$ws.Range("E12").Value = "'-1.78%"
$ws.Range("D13").Value = "'0.03540"
$ws.Range("E13").Value = "'6.68%"
$ws.Range("D14").Value = "'0.09922"
$ws.Range("E14").Value = "'-0.08%"
$ws.Range("E15").Value = "'-1.24%"
$ws.Range("D16").Value = "'0.005690"
$ws.Range("E16").Value = "'-1.39%"
$ws.Range("E17").Value = "'0.30%"
$ws.Range("D18").Value = "'2.052"
$ws.Range("E18").Value = "'-4.23%"
$ws.Range("E19").Value = "'3.08%"
$ws.Range("D20").Value = "'0.1314"
$ws.Range("E20").Value = "'-1.39%"
$ws.Range("D21").Value = "'4.555"
$ws.Range("E21").Value = "'6.28%"
$ws.Range("D22").Value = "'0.2229"
$ws.Range("E22").Value = "'-3.05%"
$ws.Range("D23").Value = "'0.04649"
$ws.Range("E23").Value = "'1.90%"
$ws.Range("D24").Value = "'0.001227"
$ws.Range("E24").Value = "'1.18%"
$ws.Range("D25").Value = "'0.004443"
$ws.Range("E25").Value = "'1.69%"
$ws.Range("D26").Value = "'0.0001298"
$ws.Range("E26").Value = "'-0.08%"
$ws.Range("D27").Value = "'0.0004744"
$ws.Range("E27").Value = "'39.73%"
$ws.Range("D39").Value = "'0.01744"
$ws.Range("E39").Value = "'-2.48%"
$ws.Range("D40").Value = "'0.04681"
$ws.Range("E40").Value = "'-2.46%"
$ws.Range("D41").Value = "'0.007845"
$ws.Range("E41").Value = "'1.14%"
$ws.Range("E42").Value = "'-1.87%"
$ws.Range("D43").Value = "'0.007657"
$ws.Range("E43").Value = "'7.75%"
$ws.Range("D44").Value = "'0.002287"
$ws.Range("E44").Value = "'7.42%"
$ws.Range("D45").Value = "'0.01007"
$ws.Range("E45").Value = "'6.69%"
$ws.Range("D46").Value = "'0.00006093"
$ws.Range("E46").Value = "'-0.26%"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("E47").Value = "'-0.08%"
$ws.Range("D48").Value = "'8.656"
$ws.Range("E48").Value = "'182.43%"
$ws.Range("E49").Value = "'35.00%"
$ws.Range("D50").Value = "'0.00002097"
$ws.Range("E50").Value = "'-0.08%"
$ws.Range("D51").Value = "'0.0001997"
$ws.Range("E51").Value = "'-0.08%"
